# try/catch for registering item
# Append the newly-registered product rows (10-19) to the inventory sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Keep the same "General" style already used by the sheet (s="1")
    # instead of letting Excel create/attach a brand new style index.
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Value = $text
}

function Set-TextCellForceText($addr, $text) {
    # For numeric-looking text (e.g. "3", "2") force it to be stored as
    # literal text (not auto-converted to a number), then restore the
    # "General" number format so the cell keeps using the shared style.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).NumberFormat = "General"
}

function Set-NumberCell($addr, $number) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Value = $number
}

try {
    # Row 10 - Glass
    Set-TextCell "A10" "1699920208662_Glass"
    Set-TextCell "B10" "Glass"
    Set-TextCellForceText "C10" "3"
    Set-TextCell "D10" "Tiny"
    Set-TextCell "E10" "For drinking"

    # Row 11 - (unnamed item)
    Set-TextCell "A11" "1699922521471_"
    Set-NumberCell "C11" 0

    # Row 12 - test
    Set-TextCell "A12" "1699923027138_test"
    Set-TextCell "B12" "test"
    Set-TextCellForceText "C12" "2"

    # Row 13 - test
    Set-TextCell "A13" "1699923083383_test"
    Set-TextCell "B13" "test"
    Set-NumberCell "C13" 0

    # Row 14 - Test
    Set-TextCell "A14" "1699923106468_Test"
    Set-TextCell "B14" "Test"
    Set-NumberCell "C14" 0

    # Row 15 - test
    Set-TextCell "A15" "1699923300586_test"
    Set-TextCell "B15" "test"
    Set-NumberCell "C15" 0

    # Row 16 - test
    Set-TextCell "A16" "1699923477814_test"
    Set-TextCell "B16" "test"
    Set-NumberCell "C16" 0

    # Row 17 - Test
    Set-TextCell "A17" "1699923506634_Test"
    Set-TextCell "B17" "Test"
    Set-NumberCell "C17" 0

    # Row 18 - Prova
    Set-TextCell "A18" "1699923598844_Prova"
    Set-TextCell "B18" "Prova"
    Set-NumberCell "C18" 0

    # Row 19 - Prova
    Set-TextCell "A19" "1699923635918_Prova"
    Set-TextCell "B19" "Prova"
    Set-NumberCell "C19" 0

    # Extend the "number stored as text" ignored-error range to cover the
    # newly registered rows (mirrors the sheet's existing A1:G9 -> A1:G19).
    try {
        $numberAsText = $ws.Range("A1:G19").Errors.Item(3)
        $numberAsText.Ignore = $true
    } catch {
        # Older/limited COM surfaces may not expose Errors.Item(xlNumberAsText);
        # harmless to skip - the data itself is already written above.
    }
}
catch {
    Write-Host "Error while registering item: $_"
}
